$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.674.31"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.793.48"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'595.18"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'166.31"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").Value = "3.791.11"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'36.25"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "4.428.54"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "3.787.02"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "'18.58"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "67.671.22"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "'10.23"
$ws.Range("E21").Value = "  -4.35%  "
$ws.Range("D22").Value = "'458.21"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'0.0000152"
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("D25").Value = "'83.72"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'11.90"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'7.30"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "'29.91"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'2.20"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "'9.21"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "3.747.03"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D44").Value = "'44.86"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").Value = "'0.299"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "'47.10"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").Value = "'8.38"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "'147.82"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").Value = "'392.44"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'1.84"
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("D51").Value = "2.756.50"
$ws.Range("E51").Value = "  +2.44%  "

Write-Output "Updated cryptos list"